$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force specific Price (column D) cells that now contain plain numeric-looking
# text (e.g. "577.53") to stay stored as TEXT, matching the source data which
# uses inline/shared strings throughout column D, not numbers.

$ws.Range("D2").Value = '66.825.10'

$ws.Range("E2").Value = '  -1.99%  '

$ws.Range("D3").Value = '3.218.49'

$ws.Range("E3").Value = '  -4.86%  '

$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '577.53'

$ws.Range("E5").Value = '  -4.79%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.07'

$ws.Range("E6").Value = '  -13.43%  '

$ws.Range("E7").Value = '  -0.02%  '

$ws.Range("D8").Value = '3.207.25'

$ws.Range("E8").Value = '  -4.88%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.525'

$ws.Range("E9").Value = '  -10.60%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.161'

$ws.Range("E10").Value = '  -14.36%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.33'

$ws.Range("E11").Value = '  -6.33%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.480'

$ws.Range("E12").Value = '  -12.18%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000234'

$ws.Range("E13").Value = '  -10.87%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.21'

$ws.Range("E14").Value = '  -15.67%  '

$ws.Range("D15").Value = '3.738.74'

$ws.Range("E15").Value = '  -5.20%  '

$ws.Range("D16").Value = '66.802.10'

$ws.Range("E16").Value = '  -2.27%  '

$ws.Range("D17").Value = '3.224.90'

$ws.Range("E17").Value = '  -5.45%  '

$ws.Range("E18").Value = '  -6.05%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.78'

$ws.Range("E19").Value = '  -14.47%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '498.44'

$ws.Range("E20").Value = '  -11.99%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.34'

$ws.Range("E21").Value = '  -14.13%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.720'

$ws.Range("E22").Value = '  -12.81%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.33'

$ws.Range("E23").Value = '  -16.25%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '82.07'

$ws.Range("E24").Value = '  -11.90%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.79'

$ws.Range("E25").Value = '  -12.44%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'

$ws.Range("E26").Value = '  -0.07%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.09'

$ws.Range("E27").Value = '  -13.26%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.03'

$ws.Range("E28").Value = '  -12.94%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '27.72'

$ws.Range("E29").Value = '  -12.85%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.57'

$ws.Range("E30").Value = '  -9.16%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.15'

$ws.Range("E31").Value = '  -4.97%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.51'

$ws.Range("E32").Value = '  -7.43%  '

$ws.Range("E33").Value = '  -0.20%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.13'

$ws.Range("E34").Value = '  -19.19%  '

$ws.Range("B35").Value = 'OKB'

$ws.Range("C35").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '54.28'

$ws.Range("E35").Value = '  -2.43%  '

$ws.Range("B36").Value = 'NEARProtocol'

$ws.Range("C36").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.35'

$ws.Range("E36").Value = '  -15.86%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '488.68'

$ws.Range("E37").Value = '  -15.74%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0417'

$ws.Range("E38").Value = '  -8.12%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0813'

$ws.Range("E39").Value = '  -12.44%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.120'

$ws.Range("E40").Value = '  -12.86%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.52'

$ws.Range("E41").Value = '  -16.67%  '

$ws.Range("D42").Value = '2.836.62'

$ws.Range("E42").Value = '  -10.53%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.55'

$ws.Range("E43").Value = '  -13.25%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.251'

$ws.Range("E44").Value = '  -11.83%  '

$ws.Range("E45").Value = '  -0.08%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.06'

$ws.Range("E46").Value = '  -10.60%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '25.11'

$ws.Range("E47").Value = '  -17.75%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '121.54'

$ws.Range("E48").Value = '  -7.31%  '

$ws.Range("D49").Value = ('0.0' + ([char]0x2083) + '0527')

$ws.Range("E49").Value = '  -18.91%  '

$ws.Range("E50").Value = '  -11.52%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.13'

$ws.Range("E51").Value = '  -21.44%  '

